$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new grading column "l1" is inserted right before the existing
# "nota_iniciativa" column (U), pushing "nota_iniciativa" one column to
# the right (to V). The new "l1" column is populated with 0 for every
# data row (2-83), matching the existing "c1"/"nota_iniciativa" columns.

$ws.Columns("U").Insert()

$ws.Range("U1").Value = "l1"
$ws.Range("U2:U83").Value = 0
